$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.530.78"
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("D3").Value = "1.670.59"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'239.06"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4774"
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("D8").Value = "'0.2622"
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("D9").Value = "'0.06165"
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("D10").Value = "1.670.07"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").Value = "'0.06987"
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("D12").Value = "'14.85"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "'0.5890"
$ws.Range("E13").Value = "  -4.23%  "
$ws.Range("D14").Value = "'4.371"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "'75.24"
$ws.Range("E15").Value = "  +3.58%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "25.526.09"
$ws.Range("E18").Value = "  +2.44%  "
$ws.Range("D19").Value = "'0.000006762"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("D21").Value = "1.885.47"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").Value = "'4.439"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "'8.735"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").Value = "'5.267"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Value = "'136.61"
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").Value = "'1.390"
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("E28").Value = "  +4.03%  "
$ws.Range("D29").Value = "'104.65"
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("D31").Value = "'0.07831"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").Value = "'3.645"
$ws.Range("E32").Value = "  +2.80%  "
$ws.Range("D33").Value = "'0.9992"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'0.04251"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("D35").Value = "'2.626"
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("D36").Value = "'0.9539"
$ws.Range("E36").Value = "  +3.66%  "
$ws.Range("D37").Value = "'0.6060"
$ws.Range("E37").Value = "  +4.42%  "
$ws.Range("D38").Value = "'2.594"
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("D39").Value = "'0.8776"
$ws.Range("E39").Value = "  +7.87%  "
$ws.Range("D40").Value = "'0.9999"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").Value = "'1.865"
$ws.Range("E41").Value = "  +3.69%  "
$ws.Range("D42").Value = "'0.01478"
$ws.Range("E42").Value = "  -4.89%  "
$ws.Range("D43").Value = "'96.38"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").Value = "'0.3757"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").Value = "'4.881"
$ws.Range("E45").Value = "  +3.22%  "
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "'6.207"
$ws.Range("E47").Value = "  +2.19%  "
$ws.Range("D48").Value = "'0.05260"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").Value = "'29.93"
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("D50").Value = "'7.412"
$ws.Range("E50").Value = "  +3.76%  "
$ws.Range("E51").Value = "  +0.25%  "
